$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-05"

# Update the header label in I1 ("2022 (through 03-04)" -> "2022 (through 03-05)")
$ws.Range("I1").Value = "2022 (through 03-05)"

# Update March (row 4) 2022 count
$ws.Range("I4").Value = 30

# Update the Total row (row 14) 2022 total
$ws.Range("I14").Value = 331
